$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New entry for row 9 (day worked 2020-07-14) ---
# Copy formatting (fill/number-format pattern) from row 7, which matches the
# "odd" row style used by row 9 (styles 9/7/5/2).
$ws.Range("B7:E7").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)
$ws.Range("C9").Value = (Get-Date -Year 2020 -Month 7 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 'Diseño de vistas para la aplicación "Pantallas" en Adobe Illustrator (MenuPrincipal)'
$ws.Range("B9").Formula = '=(IF(C9="","",1+B8))'

# --- New entry for row 8 (day worked 2020-07-15) ---
# Copy formatting from row 6, which matches the "even" row style used by
# row 8 (styles 9/13/14/15).
$ws.Range("B6:E6").Copy()
$ws.Range("B8:E8").PasteSpecial(-4122)
$ws.Range("C8").Value = (Get-Date -Year 2020 -Month 7 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 'Analisis de codigo, pruebas de codigo original, toma de decisión para rehacer el sistema de "Pantallas'
$ws.Range("B8").Formula = '=(IF(C8="","",1+B7))'

# --- Update descriptions of existing entries (append " TEST") ---
$ws.Range("E5").Value = "Documentacion de 24 Procediminetos Almacenados de la base de datos ACC MEX TEST"
$ws.Range("E6").Value = "Documentacion de 20 Procediminetos Almacenados de la base de datos ACC MEX TEST"
$ws.Range("E7").Value = "Correccion de 30 Procedimientos almacenados de la base de datos ACC MEX TEST"

# --- Hours worked on day 2 (2020-07-03) corrected from 10 to 6 ---
$ws.Range("D4").Value = 6

# --- Column E widened to fit the new, longer descriptions ---
$ws.Columns("E").ColumnWidth = 87.6

# --- Selection moved to the newly active entry row ---
$ws.Range("E13").Select()
